$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 415 (shifts old rows 415-492 down to 417-494)
$ws.Rows("415:416").Insert()

$newDate = Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0

# --- New row 415 ---
$ws.Cells.Item(415, 1).Value = 9
$ws.Cells.Item(415, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = $newDate
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = "Fruta"
$ws.Cells.Item(415, 7).Value = 100102
$ws.Cells.Item(415, 8).Value = "Cítricos"
$ws.Cells.Item(415, 9).Value = 100102005
$ws.Cells.Item(415, 10).Value = "Naranja"
$ws.Cells.Item(415, 11).Value = "Lane Late"
$ws.Cells.Item(415, 12).Value = "Primera"
$ws.Cells.Item(415, 13).Value = 550
$ws.Cells.Item(415, 14).Value = 6000
$ws.Cells.Item(415, 15).Value = 6500
$ws.Cells.Item(415, 16).Value = 6227
$ws.Cells.Item(415, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(415, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(415, 19).Value = 346
$ws.Cells.Item(415, 20).Value = 18

# --- New row 416 ---
$ws.Cells.Item(416, 1).Value = 9
$ws.Cells.Item(416, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(416, 3).Value = "Metropolitana"
$ws.Cells.Item(416, 4).Value = $newDate
$ws.Cells.Item(416, 5).Value = 13
$ws.Cells.Item(416, 6).Value = "Fruta"
$ws.Cells.Item(416, 7).Value = 100102
$ws.Cells.Item(416, 8).Value = "Cítricos"
$ws.Cells.Item(416, 9).Value = 100102005
$ws.Cells.Item(416, 10).Value = "Naranja"
$ws.Cells.Item(416, 11).Value = "New Hall"
$ws.Cells.Item(416, 12).Value = "Primera"
$ws.Cells.Item(416, 13).Value = 450
$ws.Cells.Item(416, 14).Value = 7000
$ws.Cells.Item(416, 15).Value = 7000
$ws.Cells.Item(416, 16).Value = 7000
$ws.Cells.Item(416, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(416, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(416, 19).Value = 389
$ws.Cells.Item(416, 20).Value = 18
